$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '63.959.53'
$ws.Range("E2").Value = '  -1.45%  '

$ws.Range("D3").Value = '3.149.07'
$ws.Range("E3").Value = '  -0.85%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").Value = '602.09'
$ws.Range("E5").Value = '  -2.48%  '

$ws.Range("D6").Value = '142.98'
$ws.Range("E6").Value = '  -2.78%  '

$ws.Range("E7").Value = '  +0.04%  '

$ws.Range("D8").Value = '3.142.86'
$ws.Range("E8").Value = '  -0.91%  '

$ws.Range("E9").Value = '  -0.55%  '

$ws.Range("E10").Value = '  -2.29%  '

$ws.Range("D11").Value = '5.39'
$ws.Range("E11").Value = '  -2.09%  '

$ws.Range("D12").Value = '0.467'
$ws.Range("E12").Value = '  -1.84%  '

$ws.Range("D13").Value = '0.0000254'
$ws.Range("E13").Value = '  -3.18%  '

$ws.Range("D14").Value = '35.02'
$ws.Range("E14").Value = '  -2.66%  '

$ws.Range("D15").Value = '3.663.59'
$ws.Range("E15").Value = '  -0.90%  '

$ws.Range("E16").Value = '  +2.22%  '

$ws.Range("D17").Value = '63.995.05'
$ws.Range("E17").Value = '  -1.38%  '

$ws.Range("D18").Value = '3.142.83'
$ws.Range("E18").Value = '  -0.94%  '

$ws.Range("D19").Value = '6.86'
$ws.Range("E19").Value = '  -1.36%  '

$ws.Range("D20").Value = '488.30'
$ws.Range("E20").Value = '  +1.34%  '

$ws.Range("D21").Value = '14.69'
$ws.Range("E21").Value = '  -0.63%  '

$ws.Range("E22").Value = '  -1.17%  '

$ws.Range("E23").Value = '  -2.91%  '

$ws.Range("D24").Value = '88.27'
$ws.Range("E24").Value = '  +4.15%  '

$ws.Range("D25").Value = '13.30'
$ws.Range("E25").Value = '  -4.18%  '

$ws.Range("E26").Value = '  +0.04%  '

$ws.Range("E27").Value = '  -2.49%  '

$ws.Range("D28").Value = '8.20'
$ws.Range("E28").Value = '  -5.49%  '

$ws.Range("D29").Value = '7.01'
$ws.Range("E29").Value = '  +0.83%  '

$ws.Range("E30").Value = '  -2.19%  '

$ws.Range("D31").Value = '27.73'
$ws.Range("E31").Value = '  +3.61%  '

$ws.Range("E32").Value = '  -6.03%  '

$ws.Range("E33").Value = '  +0.05%  '

$ws.Range("E34").Value = '  -2.30%  '

$ws.Range("E35").Value = '  -2.65%  '

$ws.Range("D36").Value = '6.07'
$ws.Range("E36").Value = '  +0.21%  '

$ws.Range("D38").Value = '0.0₃0748'
$ws.Range("E38").Value = '  -5.55%  '

$ws.Range("D39").Value = '2.95'
$ws.Range("E39").Value = '  -8.59%  '

$ws.Range("E40").Value = '  -1.17%  '

$ws.Range("D41").Value = '433.29'
$ws.Range("E41").Value = '  -7.56%  '

$ws.Range("E42").Value = '  -0.60%  '

$ws.Range("D43").Value = '8.37'
$ws.Range("E43").Value = '  -0.62%  '

$ws.Range("D44").Value = '2.933.62'
$ws.Range("E44").Value = '  +2.49%  '

$ws.Range("E45").Value = '  -3.53%  '

$ws.Range("E46").Value = '  -6.42%  '

$ws.Range("D47").Value = '2.40'
$ws.Range("E47").Value = '  -2.26%  '

$ws.Range("E48").Value = '  -0.10%  '

$ws.Range("D49").Value = '25.86'
$ws.Range("E49").Value = '  -3.97%  '

$ws.Range("E50").Value = '  +0.17%  '

$ws.Range("D51").Value = '120.67'
$ws.Range("E51").Value = '  -0.23%  '
